$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "23.054.12"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.590.21"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  -0.01%  "
Set-TextValue "D6" "301.73"
$ws.Range("E6").Value = "  -0.04%  "
Set-TextValue "D7" "0.3769"
$ws.Range("E7").Value = "  -0.08%  "
Set-TextValue "D8" "0.3600"
$ws.Range("E8").Value = "  -1.42%  "
Set-TextValue "D9" "50.97"
$ws.Range("E9").Value = "  +6.72%  "
Set-TextValue "D10" "1.003"
$ws.Range("E10").Value = "  +0.08%  "
Set-TextValue "D11" "1.227"
$ws.Range("E11").Value = "  -3.91%  "
Set-TextValue "D12" "0.08059"
$ws.Range("E12").Value = "  -0.20%  "
Set-TextValue "D13" "22.04"
$ws.Range("E13").Value = "  -4.02%  "
Set-TextValue "D14" "6.484"
$ws.Range("E14").Value = "  -2.22%  "
Set-TextValue "D15" "7.273"
$ws.Range("E15").Value = "  -4.82%  "
Set-TextValue "D16" "0.00001227"
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("D17").Value = "1.588.74"
$ws.Range("E17").Value = "  -0.37%  "
Set-TextValue "D18" "92.55"
$ws.Range("E18").Value = "  +1.11%  "
Set-TextValue "D19" "0.06808"
$ws.Range("E19").Value = "  +0.20%  "
Set-TextValue "D20" "17.96"
$ws.Range("E20").Value = "  -2.46%  "
Set-TextValue "D21" "6.460"
$ws.Range("E21").Value = "  -1.92%  "
Set-TextValue "D22" "1.002"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").Value = "23.058.72"
$ws.Range("E24").Value = "  -0.24%  "
Set-TextValue "D25" "2.383"
$ws.Range("E25").Value = "  +0.89%  "
Set-TextValue "D26" "2.818"
$ws.Range("E26").Value = "  -2.73%  "
Set-TextValue "D27" "20.89"
$ws.Range("E27").Value = "  -0.85%  "
Set-TextValue "D28" "148.32"
$ws.Range("E28").Value = "  -1.80%  "
Set-TextValue "D29" "5.211"
$ws.Range("E29").Value = "  -0.64%  "
Set-TextValue "D30" "133.10"
$ws.Range("E30").Value = "  +1.01%  "
Set-TextValue "D31" "2.354"
$ws.Range("E31").Value = "  -3.56%  "
Set-TextValue "D32" "6.546"
$ws.Range("E32").Value = "  -8.27%  "
$ws.Range("D33").Value = "1.767.88"
$ws.Range("E33").Value = "  -0.68%  "
Set-TextValue "D34" "0.9453"
$ws.Range("E34").Value = "  -4.10%  "
Set-TextValue "D35" "0.07397"
$ws.Range("E35").Value = "  -4.22%  "
Set-TextValue "D36" "10.16"
$ws.Range("E36").Value = "  +1.00%  "
Set-TextValue "D37" "0.02678"
$ws.Range("E37").Value = "  -3.79%  "
Set-TextValue "D38" "0.08780"
$ws.Range("E38").Value = "  -1.04%  "
Set-TextValue "D39" "6.060"
$ws.Range("E39").Value = "  -3.88%  "
Set-TextValue "D40" "0.2480"
$ws.Range("E40").Value = "  -2.59%  "
Set-TextValue "D41" "1.345"
$ws.Range("E41").Value = "  -3.68%  "
Set-TextValue "D42" "0.6927"
$ws.Range("E42").Value = "  -3.28%  "
Set-TextValue "D43" "12.11"
$ws.Range("E43").Value = "  -5.32%  "
Set-TextValue "D44" "14.98"
$ws.Range("E44").Value = "  -6.13%  "
Set-TextValue "D45" "0.6452"
$ws.Range("E45").Value = "  -2.89%  "
Set-TextValue "D46" "4.004"
$ws.Range("E46").Value = "  +0.95%  "
Set-TextValue "D47" "2.256"
$ws.Range("E47").Value = "  -2.45%  "
Set-TextValue "D48" "131.58"
$ws.Range("E48").Value = "  -0.48%  "
Set-TextValue "D49" "0.07892"
$ws.Range("E49").Value = "  -1.02%  "
Set-TextValue "D50" "1.199"
$ws.Range("E50").Value = "  +2.12%  "
Set-TextValue "D51" "1.211"
$ws.Range("E51").Value = "  +3.40%  "
